# Auto-generated Excel COM-interop edit script
# Updates cryptos list prices / 1h volume percentages, and swaps the
# Optimism / InjectiveProtocol rows (50 and 51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "598.40") must be
# forced to Text format first, otherwise Excel auto-converts them to a
# numeric cell -- the source file stores these as literal strings.
$textCells = @("D5", "D6", "D8", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D25", "D27", "D31", "D32", "D33", "D34", "D37", "D38", "D39", "D42", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '68.124.22'
$ws.Range("E2").Value = '  +1.45%  '
$ws.Range("D3").Value = '2.630.74'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '598.40'
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").Value = '153.85'
$ws.Range("E6").Value = '  +1.12%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '0.549'
$ws.Range("E8").Value = '  -1.18%  '
$ws.Range("D9").Value = '2.629.56'
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("E10").Value = '  +10.36%  '
$ws.Range("E11").Value = '  -0.62%  '
$ws.Range("D12").Value = '5.22'
$ws.Range("E12").Value = '  +0.68%  '
$ws.Range("D13").Value = '0.348'
$ws.Range("E13").Value = '  -0.37%  '
$ws.Range("D14").Value = '27.68'
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("E15").Value = '  +4.52%  '
$ws.Range("D16").Value = '3.107.65'
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").Value = '68.019.23'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").Value = '2.624.31'
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").Value = '375.10'
$ws.Range("E19").Value = '  +2.87%  '
$ws.Range("D20").Value = '11.35'
$ws.Range("E20").Value = '  +1.79%  '
$ws.Range("D21").Value = '7.48'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '4.26'
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("D23").Value = '4.83'
$ws.Range("E23").Value = '  -1.48%  '
$ws.Range("E24").Value = '  -2.53%  '
$ws.Range("D25").Value = '72.54'
$ws.Range("E25").Value = '  +7.43%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '10.00'
$ws.Range("E27").Value = '  -1.21%  '
$ws.Range("E28").Value = '  +2.98%  '
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("E30").Value = '  -1.24%  '
$ws.Range("D31").Value = '579.05'
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").Value = '1.41'
$ws.Range("E32").Value = '  +1.43%  '
$ws.Range("D33").Value = '7.86'
$ws.Range("E33").Value = '  +0.72%  '
$ws.Range("D34").Value = '1.85'
$ws.Range("E34").Value = '  +0.29%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("E36").Value = '  -0.85%  '
$ws.Range("D37").Value = '1.52'
$ws.Range("E37").Value = '  -0.85%  '
$ws.Range("D38").Value = '159.19'
$ws.Range("E38").Value = '  +0.82%  '
$ws.Range("D39").Value = '19.19'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("E40").Value = '  +5.38%  '
$ws.Range("E41").Value = '  +0.39%  '
$ws.Range("D42").Value = '5.36'
$ws.Range("E42").Value = '  +2.05%  '
$ws.Range("E43").Value = '  +2.92%  '
$ws.Range("D44").Value = '17.10'
$ws.Range("E44").Value = '  +4.55%  '
# Subscript-six (U+2086) cannot be embedded reliably via PowerShell escapes,
# build it from its code point instead.
$sub6 = [char]0x2086
$ws.Range("D45").Value = "0.0{0}0318" -f $sub6
$ws.Range("E45").Value = '  +11.80%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '40.51'
$ws.Range("E47").Value = '  -1.78%  '
$ws.Range("D48").Value = '155.45'
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("D49").Value = '3.71'
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '22.01'
$ws.Range("E50").Value = '  +7.78%  '
$ws.Range("B51").Value = 'Optimism'
$ws.Range("C51").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D51").Value = '1.70'
$ws.Range("E51").Value = '  -1.89%  '
